$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.567.25'
$ws.Range("E2").Value = '  +1.62%  '
$ws.Range("D3").Value = '1.827.81'
$ws.Range("E3").Value = '  +1.85%  '
$ws.Range("D4").Value = '0.9998'
$ws.Range("E4").Value = '  -0.10%  '
$ws.Range("D5").Value = '317.56'
$ws.Range("E5").Value = '  +0.11%  '
$ws.Range("D6").Value = '1.000'
$ws.Range("E6").Value = '  +0.01%  '
$ws.Range("D7").Value = '0.5431'
$ws.Range("E7").Value = '  +0.33%  '
$ws.Range("D8").Value = '0.4035'
$ws.Range("E8").Value = '  +6.67%  '
$ws.Range("D9").Value = '0.07678'
$ws.Range("D10").Value = '1.122'
$ws.Range("E10").Value = '  +2.53%  '
$ws.Range("D11").Value = '41.85'
$ws.Range("E11").Value = '  +0.28%  '
$ws.Range("B12").Value = 'Solana'
$ws.Range("C12").Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range("D12").Value = '21.15'
$ws.Range("E12").Value = '  +2.96%  '
$ws.Range("B13").Value = 'Polkadot'
$ws.Range("C13").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D13").Value = '6.326'
$ws.Range("E13").Value = '  +3.47%  '
$ws.Range("B14").Value = 'Chainlink'
$ws.Range("C14").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D14").Value = '7.628'
$ws.Range("E14").Value = '  +5.23%  '
$ws.Range("B15").Value = 'BinanceUSD'
$ws.Range("C15").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D15").Value = '0.9999'
$ws.Range("E15").Value = '  -0.11%  '
$ws.Range("D16").Value = '1.825.60'
$ws.Range("E16").Value = '  +2.40%  '
$ws.Range("D17").Value = '0.00001092'
$ws.Range("E17").Value = '  +3.12%  '
$ws.Range("D18").Value = '89.99'
$ws.Range("E18").Value = '  +0.98%  '
$ws.Range("D19").Value = '0.06609'
$ws.Range("E19").Value = '  +1.92%  '
$ws.Range("D20").Value = '17.80'
$ws.Range("E20").Value = '  +3.03%  '
$ws.Range("D22").Value = '6.069'
$ws.Range("E22").Value = '  +2.76%  '
$ws.Range("D23").Value = '28.568.78'
$ws.Range("E23").Value = '  +1.57%  '
$ws.Range("D24").Value = '11.17'
$ws.Range("E24").Value = '  +0.01%  '
$ws.Range("E25").Value = '  +9.24%  '
$ws.Range("D26").Value = '158.01'
$ws.Range("E26").Value = '  +1.96%  '
$ws.Range("B27").Value = 'LidoDAOToken'
$ws.Range("C27").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D27").Value = '2.457'
$ws.Range("E27").Value = '  +7.41%  '
$ws.Range("B28").Value = 'EthereumClassic'
$ws.Range("C28").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D28").Value = '20.78'
$ws.Range("E28").Value = '  +2.51%  '
$ws.Range("D29").Value = '2.036.37'
$ws.Range("E29").Value = '  +2.27%  '
$ws.Range("D30").Value = '124.09'
$ws.Range("E30").Value = '  +2.48%  '
$ws.Range("D31").Value = '1.130'
$ws.Range("E31").Value = '  +0.74%  '
$ws.Range("D32").Value = '0.1109'
$ws.Range("E32").Value = '  +5.05%  '
$ws.Range("D33").Value = '5.688'
$ws.Range("E33").Value = '  +2.39%  '
$ws.Range("D34").Value = '0.07384'
$ws.Range("E34").Value = '  +13.47%  '
$ws.Range("D35").Value = '3.645'
$ws.Range("D36").Value = '0.2247'
$ws.Range("E36").Value = '  -0.53%  '
$ws.Range("D37").Value = '0.02358'
$ws.Range("E37").Value = '  +2.79%  '
$ws.Range("E38").Value = '  +3.87%  '
$ws.Range("E39").Value = '  +5.29%  '
$ws.Range("D40").Value = '0.6302'
$ws.Range("E40").Value = '  +1.96%  '
$ws.Range("E41").Value = '  +2.48%  '
$ws.Range("D42").Value = '1.191'
$ws.Range("E42").Value = '  +1.59%  '
$ws.Range("D43").Value = '0.9994'
$ws.Range("E43").Value = '  -0.01%  '
$ws.Range("D45").Value = '13.42'
$ws.Range("E45").Value = '  +0.38%  '
$ws.Range("D46").Value = '0.5887'
$ws.Range("E46").Value = '  +1.78%  '
$ws.Range("D47").Value = '3.710'
$ws.Range("E47").Value = '  +0.93%  '
$ws.Range("D48").Value = '125.16'
$ws.Range("E48").Value = '  +0.43%  '
$ws.Range("D49").Value = '2.004'
$ws.Range("E49").Value = '  +4.16%  '
$ws.Range("D50").Value = '1.199'
$ws.Range("E50").Value = '  +0.52%  '
$ws.Range("D51").Value = '0.06910'
$ws.Range("E51").Value = '  +1.48%  '
